$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.945.09'

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.363.00'
$ws.Range("E3").Value = '  +2.17%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.03%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '302.60'
$ws.Range("E5").Value = '  +0.31%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '95.84'
$ws.Range("E6").Value = '  +0.44%  '

# Row 7
$ws.Range("E7").Value = '  -0.03%  '

# Row 9
$ws.Range("E9").Value = '  -0.37%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.14'
$ws.Range("E10").Value = '  -0.04%  '

# Row 11
$ws.Range("E11").Value = '  +3.89%  '

# Row 12
$ws.Range("E12").Value = '  +0.21%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.31'
$ws.Range("E13").Value = '  -3.41%  '

# Row 14
$ws.Range("E14").Value = '  -0.11%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.729.98'
$ws.Range("E15").Value = '  +2.14%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.380.12'
$ws.Range("E16").Value = '  +2.71%  '

# Row 17
$ws.Range("E17").Value = '  +0.54%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.909.13'
$ws.Range("E18").Value = '  +0.48%  '

# Row 19
$ws.Range("B19").Value = 'Uniswap'
$ws.Range("C19").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.26'
$ws.Range("E19").Value = '  +2.04%  '

# Row 20
$ws.Range("B20").Value = 'InternetComputer(DFINITY)'
$ws.Range("C20").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.85'
$ws.Range("E20").Value = '  -2.28%  '

# Row 21
$ws.Range("E21").Value = '  -0.84%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.92'
$ws.Range("E22").Value = '  +0.28%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '234.98'
$ws.Range("E23").Value = '  -0.06%  '

# Row 24
$ws.Range("E24").Value = '  -5.25%  '

# Row 25
$ws.Range("E25").Value = '  -0.01%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.42'
$ws.Range("E26").Value = '  +0.44%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.46'
$ws.Range("E27").Value = '  +0.59%  '

# Row 28
$ws.Range("E28").Value = '  +0.80%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.31'
$ws.Range("E29").Value = '  +2.20%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '32.06'
$ws.Range("E30").Value = '  -0.67%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'

# Row 32
$ws.Range("E32").Value = '  +0.29%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '17.34'
$ws.Range("E33").Value = '  -2.15%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0714'
$ws.Range("E34").Value = '  +2.29%  '

# Row 35
$ws.Range("B35").Value = 'Monero'
$ws.Range("C35").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '128.32'
$ws.Range("E35").Value = '  -22.64%  '

# Row 36
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.85'
$ws.Range("E36").Value = '  +3.23%  '

# Row 37
$ws.Range("E37").Value = '  +3.34%  '

# Row 38
$ws.Range("E38").Value = '  -2.33%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.28'
$ws.Range("E39").Value = '  -2.44%  '

# Row 40
$ws.Range("E40").Value = '  +2.60%  '

# Row 41
$ws.Range("E41").Value = '  -0.76%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '21.21'
$ws.Range("E42").Value = '  -2.16%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.929.94'
$ws.Range("E43").Value = '  +0.40%  '

# Row 44
$ws.Range("E44").Value = '  -0.32%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.13'
$ws.Range("E45").Value = '  +2.57%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.71'
$ws.Range("E46").Value = '  -0.88%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.14'
$ws.Range("E47").Value = '  -9.03%  '

# Row 48
$ws.Range("B48").Value = 'Stacks'
$ws.Range("C48").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.50'
$ws.Range("E48").Value = '  +1.79%  '

# Row 49
$ws.Range("B49").Value = 'MultiversX'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '51.50'
$ws.Range("E49").Value = '  -3.31%  '

# Row 50
$ws.Range("E50").Value = '  +1.01%  '

# Row 51
$ws.Range("B51").Value = 'BitcoinSV'
$ws.Range("C51").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '71.07'
$ws.Range("E51").Value = '  -1.51%  '
